$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title cell B1 (project/commission code changed)
$ws.Range("B1").Value = "CEMRJ1CO_1013"

# Apply the bordered cell format (same look as other data cells) to the
# new block before filling in values, mirroring a copy-format + remove
# wrap-text operation.
$ws.Range("A1").Copy()
$ws.Range("A4:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A4:G6").WrapText = $false

# Add new data rows 4-6 (dossier IXISPI)
$ws.Range("A4").Value = "IMB/26198/C/03EG"
$ws.Range("B4").Value = "CEMRJ1CO_1013"
$ws.Range("C4").Value = "Obtention BPT"
$ws.Range("D4").Value = "ZN_0319_26_0075"
$ws.Range("E4").Value = "SO"
$ws.Range("F4").Value = "SO"
$ws.Range("G4").Value = "En attente BPT"

$ws.Range("A5").Value = "IMB/26198/C/03EG"
$ws.Range("B5").Value = "CEMRJ1CO_1013"
$ws.Range("C5").Value = "Obtention BPT"
$ws.Range("D5").Value = "SO"
$ws.Range("E5").Value = "ZE_1903_26_0017_01"
$ws.Range("F5").Value = "RGT_1903_26_0046"
$ws.Range("G5").Value = "Qualif en cours"

$ws.Range("A6").Value = "IMB/26198/C/03ZV"
$ws.Range("B6").Value = "CEMRJ1CO_1013"
$ws.Range("C6").Value = "Obtention BPT"
$ws.Range("D6").Value = "ZN_0319_26_0083"
$ws.Range("E6").Value = "SO"
$ws.Range("F6").Value = "SO"
$ws.Range("G6").Value = "En attente BPT"
